$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Helper: appends "  Can be specified using ${Property}." to the end of the
# Description cell (column 2) of the given table row, with the "${Property}"
# portion formatted using the RTiSWDocLiteralText character style, matching
# the style used elsewhere in this document for literal/code text.
function Add-CanBeSpecifiedUsingProperty($RowIndex) {

    $cell = $t.Cell($RowIndex, 2)
    $cellRange = $cell.Range

    # Position right before the cell end-of-cell mark (last char of Range).
    $insertPoint = $cellRange.End - 1

    # Insert the leading plain text.
    $leadRange = $d.Range($insertPoint, $insertPoint)
    $leadRange.InsertAfter('  Can be specified using ')

    # Insert the literal "${Property}" text right after the lead text.
    $propStart = $insertPoint + '  Can be specified using '.Length
    $propRange = $d.Range($propStart, $propStart)
    $propRange.InsertAfter('${Property}')

    # Insert trailing period.
    $trailStart = $propStart + '${Property}'.Length
    $trailRange = $d.Range($trailStart, $trailStart)
    $trailRange.InsertAfter('.')

    # Apply the RTiSWDocLiteralText character style to just "${Property}".
    # Using Range.Style directly does not reliably produce a character run
    # style inside table cells with this runtime, so use Find/Replace with
    # a Replacement.Style instead, scoped to the small range we just wrote.
    $styleTargetStart = $propStart
    $styleTargetEnd = $propStart + '${Property}'.Length
    $styleRange = $d.Range($styleTargetStart, $styleTargetEnd)

    $find = $styleRange.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Replacement.Style = "RTiSWDocLiteralText"
    # Wrap = 0 (wdFindStop) so the search/replace stays confined to the
    # small scoped $styleRange instead of leaking into the rest of the document.
    $find.Execute('${Property}', $false, $false, $false, $false, $false, $true, 0, $false, '${Property}', 2) | Out-Null
}

# Row 3 = TSID
Add-CanBeSpecifiedUsingProperty 3
# Row 4 = EnsembleID
Add-CanBeSpecifiedUsingProperty 4
# Row 7 = LongitudeProperty
Add-CanBeSpecifiedUsingProperty 7
# Row 8 = LatitudeProperty
Add-CanBeSpecifiedUsingProperty 8
# Row 9 = ElevationProperty
Add-CanBeSpecifiedUsingProperty 9
# Row 10 = WKTGeometryProperty
Add-CanBeSpecifiedUsingProperty 10
# Row 11 = IncludeColumns
Add-CanBeSpecifiedUsingProperty 11
# Row 12 = ExcludeColumns
Add-CanBeSpecifiedUsingProperty 12
# Row 13 = JavaScriptVar
Add-CanBeSpecifiedUsingProperty 13
